# Update the marksheet's "Marking" (per-correct-answer score) and the
# derived "Total" / fraction-of-max cells to reflect the new marking scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row (row 11): points awarded per correct answer, 3 -> 5
$ws.Range("B11").Value = 5

# Total row (row 12): recomputed total correct marks, 60 -> 100
$ws.Range("B12").Value = 100

# Total row (row 12): fraction "correct/max" text, "54/84" -> "100/140"
$ws.Range("E12").Value = "100/140"
